# This script re-applies an upstream "automatic update of files" sync: the
# records (observations) shown in rows 9-14 and 16-18 were refreshed from the
# source system, which re-ordered/re-matched several sightings of the same
# species onto different coordinate points. Net effect on the worksheet is
# that the "record" values in columns A (Id), B (Taxonsorteringsordning),
# D (Rodlistade), E (TaxonId), F (Artnamn), G (Vetenskapligt namn),
# H (Auktor), Q (Ost) and R (Nord) move between rows while every other
# column (location names, dates, observer, etc.) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RecordRow {
    param($Row, $A, $B, $D, $E, $F, $G, $H, $Q, $R)

    $ws.Cells.Item($Row, 1).Value  = $A   # A: Id
    $ws.Cells.Item($Row, 2).Value  = $B   # B: Taxonsorteringsordning
    $ws.Cells.Item($Row, 4).Value  = $D   # D: Rodlistade
    $ws.Cells.Item($Row, 5).Value  = $E   # E: TaxonId
    $ws.Cells.Item($Row, 6).Value  = $F   # F: Artnamn
    $ws.Cells.Item($Row, 7).Value  = $G   # G: Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value  = $H   # H: Auktor
    $ws.Cells.Item($Row, 17).Value = $Q   # Q: Ost
    $ws.Cells.Item($Row, 18).Value = $R   # R: Nord
}

# New values for each row, taken from the refreshed source record that now
# belongs there.
Set-RecordRow 9  111670599 96348 "VU" 220787 "Knärot"          "Goodyera repens"    "(L.) R. Br." 558031.5226908802 7067909.315233406
Set-RecordRow 10 111671395 96348 "VU" 220787 "Knärot"          "Goodyera repens"    "(L.) R. Br." 557763.2623863788 7068264.582601988
Set-RecordRow 11 111670575 96346 "NT" 620    "Skogsfru"        "Epipogium aphyllum" "Sw."         558082.6649719321 7067974.943554637
Set-RecordRow 12 111670588 96348 "VU" 220787 "Knärot"          "Goodyera repens"    "(L.) R. Br." 558039.6361001397 7067902.375451046
Set-RecordRow 13 111671345 96348 "VU" 220787 "Knärot"          "Goodyera repens"    "(L.) R. Br." 557812.5300353739 7068166.248475613
Set-RecordRow 14 111671364 96368 "LC" 221952 "Spindelblomster" "Neottia cordata"    "(L.) Rich."  557813.3601359134 7068169.364891288
Set-RecordRow 16 111670593 78578 "NT" 6458   "Lunglav"         "Lobaria pulmonaria" "(L.) Hoffm." 558040.5475534229 7067901.063021242
Set-RecordRow 17 111671384 96348 "VU" 220787 "Knärot"          "Goodyera repens"    "(L.) R. Br." 557798.0632258818 7068181.046264404
Set-RecordRow 18 111670607 96368 "LC" 221952 "Spindelblomster" "Neottia cordata"    "(L.) Rich."  558031.5471372061 7067907.98648507

# Column L ("Kön") is an otherwise-blank marker column. It previously had an
# (empty) placeholder cell on row 16 and none on row 17; after the refresh
# above, that placeholder now belongs with the record that landed on row 17
# instead of row 16.
$ws.Cells.Item(16, 12).ClearContents()
$ws.Cells.Item(17, 12).Value = ""
